$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells being updated, so numeric-looking
# strings (e.g. "1.00", "308.51") are preserved as text, matching source data type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.676.70'
$ws.Range("D3").Value = '2.534.35'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").Value = '308.51'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").Value = '100.57'
$ws.Range("E6").Value = '  +3.56%  '
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  -2.49%  '
$ws.Range("D10").Value = '36.14'
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("D11").Value = '0.0804'
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").Value = '2.930.45'
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").Value = '15.90'
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("D16").Value = '2.561.68'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("E17").Value = '  -4.29%  '
$ws.Range("D18").Value = '42.648.24'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").Value = '12.23'
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").Value = '69.36'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '243.85'
$ws.Range("E23").Value = '  -4.11%  '
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '26.00'
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("D28").Value = '2.32'
$ws.Range("E28").Value = '  -4.94%  '
$ws.Range("D29").Value = '39.19'
$ws.Range("E29").Value = '  -3.03%  '
$ws.Range("D30").Value = '10.17'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").Value = '155.56'
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  +12.91%  '
$ws.Range("D34").Value = '0.0792'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("E35").Value = '  -2.54%  '
$ws.Range("D36").Value = '18.39'
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("D37").Value = '2.03'
$ws.Range("E37").Value = '  -5.55%  '
$ws.Range("E38").Value = '  -7.19%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  +7.78%  '
$ws.Range("D42").Value = '22.01'
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '3.29'
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").Value = '0.0298'
$ws.Range("E45").Value = '  -2.11%  '
$ws.Range("D46").Value = '1.966.79'
$ws.Range("E46").Value = '  -1.91%  '
$ws.Range("D47").Value = '8.89'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").Value = '80.70'
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.192'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.727.32'
$ws.Range("E50").Value = '  -3.56%  '
$ws.Range("D51").Value = '0.850'
$ws.Range("E51").Value = '  +8.87%  '
